{"js": "// Update the lattice-multiplication practice table: every cell keeps its\n// existing 5-line layout (problem header, split digits, the \"----\"\n// separator, and the two lattice working rows) but the numbers inside\n// are replaced with a new set of practice problems. The table's shape\n// (5 rows x 3 columns x 1 paragraph/run per cell) is unchanged; only the\n// <w:t> text content changes.\n//\n// New 5-line content for every cell, in row-major order (5 rows x 3 cols),\n// taken line-for-line from the target diff.\nconst newCells = [\n  [\n    [\"66 x 36\", \"  3    6\", \"  ----\", \"6|    |\", \"6|    |\"],\n    [\"94 x 94\", \"  9    4\", \"  ----\", \"9|    |\", \"4|    |\"],\n    [\"30 x 73\", \"  7    3\", \"  ----\", \"3|    |\", \"0|    |\"],\n  ],\n  [\n    [\"91 x 39\", \"  3    9\", \"  ----\", \"9|    |\", \"1|    |\"],\n    [\"42 x 74\", \"  7    4\", \"  ----\", \"4|    |\", \"2|    |\"],\n    [\"69 x 10\", \"  1    0\", \"  ----\", \"6|    |\", \"9|    |\"],\n  ],\n  [\n    [\"84 x 13\", \"  1    3\", \"  ----\", \"8|    |\", \"4|    |\"],\n    [\"54 x 52\", \"  5    2\", \"  ----\", \"5|    |\", \"4|    |\"],\n    [\"39 x 65\", \"  6    5\", \"  ----\", \"3|    |\", \"9|    |\"],\n  ],\n  [\n    [\"58 x 15\", \"  1    5\", \"  ----\", \"5|    |\", \"8|    |\"],\n    [\"32 x 45\", \"  4    5\", \"  ----\", \"3|    |\", \"2|    |\"],\n    [\"37 x 83\", \"  8    3\", \"  ----\", \"3|    |\", \"7|    |\"],\n  ],\n  [\n    [\"84 x 20\", \"  2    0\", \"  ----\", \"8|    |\", \"4|    |\"],\n    [\"88 x 51\", \"  5    1\", \"  ----\", \"8|    |\", \"8|    |\"],\n    [\"99 x 87\", \"  8    7\", \"  ----\", \"9|    |\", \"9|    |\"],\n  ],\n];\n\nfunction escapeXml(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n}\n\n// Build a <w:p> for a cell from its 5 lines, joined by <w:br/>, reproducing\n// the source formatting: a single run at sz=32, with xml:space=\"preserve\"\n// on any line that has leading/trailing spaces (so the whitespace survives\n// round-tripping exactly like the original markup).\nfunction cellParagraphXml(lines) {\n  const runsXml = lines\n    .map((line, i) => {\n      const needsPreserve = /^\\s|\\s$/.test(line);\n      const spaceAttr = needsPreserve ? ' xml:space=\"preserve\"' : \"\";\n      const tEl = `<w:t${spaceAttr}>${escapeXml(line)}</w:t>`;\n      return i === 0 ? tEl : `<w:br/>${tEl}`;\n    })\n    .join(\"\");\n  return `<w:p><w:r><w:rPr><w:sz w:val=\"32\"/></w:rPr>${runsXml}</w:r></w:p>`;\n}\n\n// Office.js's Range.insertOoxml requires a full \"flat OPC\" package wrapper\n// around the OOXML fragment being inserted.\nfunction flatOpcDocument(paragraphXml) {\n  return (\n    '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body>\" +\n    paragraphXml +\n    \"</w:body></w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nfor (let r = 0; r < rows.items.length; r++) {\n  const cells = rows.items[r].cells.items;\n  for (let c = 0; c < cells.length; c++) {\n    const xml = flatOpcDocument(cellParagraphXml(newCells[r][c]));\n    cells[c].body.insertOoxml(xml, \"Replace\");\n  }\n}\nawait context.sync();\n", "ps1": "# Update the lattice-multiplication practice table: every cell keeps its\n# existing 5-line layout (problem header, split digits, the \"----\"\n# separator, and the two lattice working rows) but the numbers inside are\n# replaced with a new set of practice problems. The table's shape\n# (5 rows x 3 columns) is unchanged; only each cell's text content changes.\n#\n# [char]11 is a vertical-tab, which is how Word represents a manual line\n# break (<w:br/>) inside Range.Text.\n$LF = [char]11\n\n# New 5-line content for every cell, in row-major order (5 rows x 3 cols),\n# taken line-for-line from the target diff.\n$newCells = @(\n  @(\n    @(\"66 x 36\", \"  3    6\", \"  ----\", \"6|    |\", \"6|    |\"),\n    @(\"94 x 94\", \"  9    4\", \"  ----\", \"9|    |\", \"4|    |\"),\n    @(\"30 x 73\", \"  7    3\", \"  ----\", \"3|    |\", \"0|    |\")\n  ),\n  @(\n    @(\"91 x 39\", \"  3    9\", \"  ----\", \"9|    |\", \"1|    |\"),\n    @(\"42 x 74\", \"  7    4\", \"  ----\", \"4|    |\", \"2|    |\"),\n    @(\"69 x 10\", \"  1    0\", \"  ----\", \"6|    |\", \"9|    |\")\n  ),\n  @(\n    @(\"84 x 13\", \"  1    3\", \"  ----\", \"8|    |\", \"4|    |\"),\n    @(\"54 x 52\", \"  5    2\", \"  ----\", \"5|    |\", \"4|    |\"),\n    @(\"39 x 65\", \"  6    5\", \"  ----\", \"3|    |\", \"9|    |\")\n  ),\n  @(\n    @(\"58 x 15\", \"  1    5\", \"  ----\", \"5|    |\", \"8|    |\"),\n    @(\"32 x 45\", \"  4    5\", \"  ----\", \"3|    |\", \"2|    |\"),\n    @(\"37 x 83\", \"  8    3\", \"  ----\", \"3|    |\", \"7|    |\")\n  ),\n  @(\n    @(\"84 x 20\", \"  2    0\", \"  ----\", \"8|    |\", \"4|    |\"),\n    @(\"88 x 51\", \"  5    1\", \"  ----\", \"8|    |\", \"8|    |\"),\n    @(\"99 x 87\", \"  8    7\", \"  ----\", \"9|    |\", \"9|    |\")\n  )\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n  for ($c = 1; $c -le $t.Columns.Count; $c++) {\n    $lines = $newCells[$r - 1][$c - 1]\n    $text = [string]::Join($LF, $lines)\n    $cell = $t.Cell($r, $c)\n    $cell.Range.Text = $text\n  }\n}\n"}
